# Deliverables Tracking workbook update
# - Adds a new "Make update to Architecture" deliverable row to the
#   "Hardware Development Process" sheet (Realization section), and
# - Fills in the Class Assigned / Assigned / Due Date / Time columns for
#   the rest of that Realization group (Estimate Task Hours, Eagle
#   Library, Eagle Schematic, Eagle Layout, Cable Definition), which had
#   previously been left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hardware Development Process")
$ws.Activate()

# Insert a new row above the existing "Estimate Task Hours" row (row 17)
# to make room for the new "Make update to Architecture" deliverable.
# Excel copies the row-above's formatting onto the freshly inserted row.
$ws.Rows.Item(17).Insert()

# New deliverable row.
$ws.Range("A17").Value = "Make update to Architecture"

# Copy the existing date formatting (d-mmm, used elsewhere in this
# column) onto the newly touched date cells before writing values, so
# the new cells line up with the rest of the sheet instead of getting a
# generic date format.
$ws.Range("D6:E6").Copy()
$ws.Range("D17:E22").PasteSpecial(-4122)

$ws.Range("C17:C22").Value = 3
$ws.Range("D17:D22").Value = 42879

$ws.Range("E17").Value = 42885
$ws.Range("E18").Value = 42912
$ws.Range("E19").Value = 42912
$ws.Range("E20").Value = 42912
$ws.Range("E21").Value = 42912
$ws.Range("E22").Value = 42912

$ws.Range("F17:F22").Value = "End of Day"

# Leave the selection where data entry finished.
$ws.Range("F18").Select()
